$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.127881588408715
$ws.Range("C2").Value = 114.8270160096505
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 8.660232485948974
$ws.Range("G2").Value = 124.4206197205922

# Row 3
$ws.Range("B3").Value = 0.01514828764759746
$ws.Range("C3").Value = 0.3127903958511391
$ws.Range("D3").Value = 0.1575252929769615
$ws.Range("E3").Value = 0.496779210170732
$ws.Range("G3").Value = 0.9822431866464301

# Row 4
$ws.Range("B4").Value = 3.230985683306322
$ws.Range("C4").Value = 1.667794583268128
$ws.Range("D4").Value = 3.900430680208489
$ws.Range("E4").Value = 0.496779210170732
$ws.Range("G4").Value = 9.295990156953671

# Row 5
$ws.Range("B5").Value = 3.230985683306322
$ws.Range("C5").Value = 1.667794583268128
$ws.Range("D5").Value = 0.8054896365839992
$ws.Range("E5").Value = 0.496779210170732
$ws.Range("G5").Value = 6.201049113329182
